# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values computed for rows 2..19 (column G), replacing the previous Strike# values.
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 3
    16 = 2
    17 = 0
    18 = 1
    19 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
